$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projects")

$ws.Range("A3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("C6").ClearContents()

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F11").Select()
